$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. DTR summary: drop the stray SICK LEAVE (I) entries that had crept into
#    the NO. OF HOURS UNDERTIME (F) column's rows, moving their values into F.
# ---------------------------------------------------------------------------
$ws.Range("F8").Value = 2.75
$ws.Range("I8").Value = ""

$ws.Range("F9").Value = 2.0

$ws.Range("I10").Value = ""

$ws.Range("I16").Value = ""

# ---------------------------------------------------------------------------
# 2. Add a "Legends:" section below the summary table (rows 24-30).
# ---------------------------------------------------------------------------

# "Legends:" heading -- reuse the same bold/underline/size-15 look as the
# report title (A1) by copying its format.
$ws.Range("E24").Value = "Legends:"
$ws.Range("A1").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("E24:P24").Merge()

# Row 25-26: blue swatch + explanation (request/remark day).
$ws.Range("E25").Interior.Color = 13411113
$ws.Range("F25").Value = "Employee has request(s)/remark(s) for that day.`n*May incur late and/or undertime depending on his or her time-in and time-out."
$ws.Range("F25").Font.Bold = $true
$ws.Range("F25").Font.Underline = $true
$ws.Range("E25:E26").Merge()
$ws.Range("F25:P26").Merge()

# Row 27-28: orange swatch + explanation (half-day).
$ws.Range("E27").Interior.Color = 6737151
$ws.Range("F27").Value = "Employee is considered half-day because of his time-in or time-out."
$ws.Range("F27").Font.Bold = $true
$ws.Range("F27").Font.Underline = $true
$ws.Range("E27:E28").Merge()
$ws.Range("F27:P28").Merge()

# Row 29-30: red swatch + explanation (absent).
$ws.Range("E29").Interior.Color = 6184671
$ws.Range("F29").Value = "Employee has no time-in and therefore, considered as absent."
$ws.Range("F29").Font.Bold = $true
$ws.Range("F29").Font.Underline = $true
$ws.Range("E29:E30").Merge()
$ws.Range("F29:P30").Merge()
